# Planification.xlsx edit script
# Implements:
#  - SPRINT 0 / SPRINT 1 / SPRINT 2: add the A1 "total hours" formula
#  - SPRINT 2: add row 9 (a newly logged task)
#  - New sheet "SPRINT 3" with its own task rows + A1 total formula
#  - SPRINTS: insert a "totals / deadline" row at the top (shifting existing
#    rows down) + a deadline block at the bottom, new formulas, new columns
#  - AGENDA: view/selection tweak
#  - Misc sheet selections / active tab bookkeeping
#
# NOTE on ordering: new text values are written in the exact sequence the
# original author entered them so that freshly-created shared-string indices
# line up (65 .. 78) the same way they do in the authored workbook.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$sprint0 = $wb.Worksheets.Item("SPRINT 0")
$sprint1 = $wb.Worksheets.Item("SPRINT 1")
$sprint2 = $wb.Worksheets.Item("SPRINT 2")
$sprints = $wb.Worksheets.Item("SPRINTS")
$agenda  = $wb.Worksheets.Item("AGENDA")

# ---------------------------------------------------------------------------
# 1) SPRINT 2 (existing sheet) - add a new logged task (row 9)
# ---------------------------------------------------------------------------
$sprint2.Range("B8").Copy()
$sprint2.Range("B9").PasteSpecial($xlPasteFormats)
$sprint2.Range("B9").Value = 43930
$sprint2.Range("C9").Value = "Planification, sprint review, Copyleaks test"
$sprint2.Range("D9").Value = 5

$sprint2.Range("D17").Formula = "=SUM(D3:D11)"

# ---------------------------------------------------------------------------
# 2) New sheet "SPRINT 3", positioned right after "SPRINT 2"
# ---------------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Add($null, $sprint2)
$sprint3.Name = "SPRINT 3"

# Match the formatting of an existing SPRINT sheet (column width, header
# row, date-formatted column, borders on the blank trailing rows, etc.)
$sprint2.Columns.Item(3).Copy()
$sprint3.Columns.Item(3).PasteSpecial($xlPasteFormats)

$sprint2.Range("B2:D17").Copy()
$sprint3.Range("B2:D17").PasteSpecial($xlPasteFormats)

$sprint3.Range("B2").Value = "Date"
$sprint3.Range("C2").Value = "Quoi"
$sprint3.Range("D2").Value = "Temps (h)"

$sprint3.Range("B3").Value = 43932
$sprint3.Range("C3").Value = "Unicheck interface web"
$sprint3.Range("D3").Value = 2

$sprint3.Range("B4").Value = 43937
$sprint3.Range("C4").Value = "Unicheck test API"
$sprint3.Range("D4").Value = 3

$sprint3.Range("B5").Value = 43937
$sprint3.Range("C5").Value = "Unicheck appel support technique, email"
$sprint3.Range("D5").Value = 1

$sprint3.Range("B6").Value = 43937
$sprint3.Range("C6").Value = "Planification"
$sprint3.Range("D6").Value = 0.5

$sprint3.Range("B7").Value = 43937
$sprint3.Range("C7").Value = "Test connections FTP, HTTP"
$sprint3.Range("D7").Value = 0.5

$sprint3.Range("B8").Value = 43938
$sprint3.Range("C8").Value = "Test MYSQL connection"
$sprint3.Range("D8").Value = 0.25

$sprint3.Range("B9").Value = 43938
$sprint3.Range("C9").Value = "Dessiner schéma collaboratif du projet"
$sprint3.Range("D9").Value = 0.5

$sprint3.Range("B10").Value = 43938
$sprint3.Range("C10").Value = "Test page local XAMPP"
$sprint3.Range("D10").Value = 0.5

$sprint3.Range("B11").ClearContents()

$sprint3.Range("D17").Formula = "=SUM(D3:D11)"

$sprint3.Range("A2").Select()

# ---------------------------------------------------------------------------
# 3) SPRINT 0 / SPRINT 1 / SPRINT 2 - the A1 "hours logged so far" formula
# ---------------------------------------------------------------------------
$sprint0.Range("A1").Formula = "=SUM(D3:D14)"
$sprint1.Range("A1").Formula = "=SUM(D3:D11)"
$sprint2.Range("A1").Formula = "=SUM(D3:D11)"
$sprint3.Range("A1").Formula = "=SUM(D3:D11)"

# ---------------------------------------------------------------------------
# 4) SPRINTS sheet - insert a new totals row at the top, shifting the rest
#    down by one, then append the deadline block at the bottom.
# ---------------------------------------------------------------------------
$sprints.Rows.Item(1).Insert()

$sprints.Range("A1").Value = "Total projet"
$sprints.Range("B1").Formula = "='SPRINT 0'!A1+'SPRINT 1'!A1+'SPRINT 2'!A1+'SPRINT 3'!A1"

$sprints.Range("A12").Value = "Temps avant rendu"

$sprints.Range("B11").Value = "Deadline"

$sprints.Range("E1").Value = "Temps total à réaliser"

$sprints.Range("D11").Value = "Cmb travailler par j"

$sprint0.Range("B3").Copy()
$sprints.Range("C11").PasteSpecial($xlPasteFormats)
$sprints.Range("C11").Value = 44043

$sprints.Range("C12").Formula = "=C11-TODAY()"
$sprints.Range("C12").NumberFormat = "0.00"

$sprints.Range("D12").Formula = "=(F1-B1)/C12"

$sprints.Range("F1").Value = 330

$sprints.Columns.Item(3).ColumnWidth = 8.830729166666666
$sprints.Columns.Item(5).ColumnWidth = 16.963541666666668

$sprints.PageSetup.Orientation = 1
$sprints.PageSetup.PaperSize = 9

# ---------------------------------------------------------------------------
# 5) AGENDA sheet - clear the scroll position, move the selection
# ---------------------------------------------------------------------------
$agenda.Activate()
$agenda.Range("A3").Select()
$agenda.Range("J11").Select()

# ---------------------------------------------------------------------------
# 6) Reset selections on the other sheets to A2 (matches authored state)
# ---------------------------------------------------------------------------
$sprint0.Activate()
$sprint0.Range("A2").Select()

$sprint1.Activate()
$sprint1.Range("A2").Select()

$sprint2.Activate()
$sprint2.Range("A2").Select()

$sprint3.Activate()
$sprint3.Range("A2").Select()

# ---------------------------------------------------------------------------
# 7) Make SPRINTS the active tab / sheet, with its authored selection
# ---------------------------------------------------------------------------
$sprints.Activate()
$sprints.Range("F15").Select()
